# The workbook originally had two sheets, each containing the same
# "nutrient density" table (columns A:H) plus a second, duplicate copy of
# the same table pasted into columns I:O of sheet 1 (mirroring sheet 2's
# data). This edit:
#   1. Renames the sheets to shorter titles.
#   2. Removes the duplicate block (columns I:O) from sheet 1, shrinking
#      its used range back down to A1:H36.
#   3. Makes sheet 1 the active/selected tab (it was sheet 2 before),
#      carrying over the leftover stale selection on sheet 1 and giving
#      sheet 2 a fresh selection.

$wb  = $excel.ActiveWorkbook
$ws1 = $wb.Worksheets.Item(1)
$ws2 = $wb.Worksheets.Item(2)

# 1. Rename the sheet tabs.
$ws1.Name = "2015 2016"
$ws2.Name = "2017 2018"

# 2. Delete the duplicated columns I:O from the first sheet.
$ws1.Columns("I:O").Delete()

# 3. Update the active sheet / selections.
$null = $ws2.Range("N22").Select()
$ws1.Activate()
